$d = $word.ActiveDocument

$ids = @("p002v_1", "p002v_2", "p002v_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $r = $d.Content
    $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $old, 2)
}
